$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.562782333333333
$ws.Range("H2").Value = 7.688347
$ws.Range("I2").Value = 0.9132775872694281
$ws.Range("J2").Value = 0.9404640236759205
$ws.Range("M2").Value = 0.4260053333333333
$ws.Range("N2").Value = 1.278016
$ws.Range("O2").Value = 0.02405532912416773
$ws.Range("P2").Value = 0.02531756756689831
$ws.Range("Q2").Value = 1.091758942172445
$ws.Range("R2").Value = 9.825830479552002
$ws.Range("S2").Value = 0.0219691929434919
$ws.Range("T2").Value = 0.02381026146365217

$ws.Range("G3").Value = 2.562782333333333
$ws.Range("H3").Value = 7.688347
$ws.Range("I3").Value = 0.9132775872694281
$ws.Range("J3").Value = 0.9404640236759205
$ws.Range("O3").Value = 0.05879323641880037
$ws.Range("P3").Value = 0.06187825274916518
$ws.Range("Q3").Value = 2.668350171729555
$ws.Range("R3").Value = 24.015151545566
$ws.Range("S3").Value = 0.05369454510432307
$ws.Range("T3").Value = 0.05819427055851548

$ws.Range("G4").Value = 2.562782333333333
$ws.Range("H4").Value = 7.688347
$ws.Range("I4").Value = 0.9132775872694281
$ws.Range("J4").Value = 0.9404640236759205
$ws.Range("M4").Value = 5.850740666666667
$ws.Range("N4").Value = 17.552222
$ws.Range("O4").Value = 0.3303749538898241
$ws.Range("P4").Value = 0.3477104875323931
$ws.Range("Q4").Value = 14.99417481744822
$ws.Range("R4").Value = 134.947573357034
$ws.Range("S4").Value = 0.3017240407827471
$ws.Range("T4").Value = 0.3270092041790305

$ws.Range("G5").Value = 2.562782333333333
$ws.Range("H5").Value = 7.688347
$ws.Range("I5").Value = 0.9132775872694281
$ws.Range("J5").Value = 0.9404640236759205
$ws.Range("M5").Value = 2.648771
$ws.Range("N5").Value = 5.297542
$ws.Range("O5").Value = 0.1495686865725097
$ws.Range("P5").Value = 0.1049445996947469
$ws.Range("Q5").Value = 6.788223523845667
$ws.Range("R5").Value = 40.729341143074
$ws.Range("S5").Value = 0.136597729203999
$ws.Range("T5").Value = 0.09869662049198041

$ws.Range("G6").Value = 2.562782333333333
$ws.Range("H6").Value = 7.688347
$ws.Range("I6").Value = 0.9132775872694281
$ws.Range("J6").Value = 0.9404640236759205
$ws.Range("M6").Value = 7.742685666666667
$ws.Range("N6").Value = 23.228057
$ws.Range("O6").Value = 0.4372077939946981
$ws.Range("P6").Value = 0.4601490924567965
$ws.Range("Q6").Value = 19.84281803908656
$ws.Range("R6").Value = 178.585362351779
$ws.Range("S6").Value = 0.399292079234867
$ws.Range("T6").Value = 0.432753666982742

$ws.Range("G7").Value = 0.243355
$ws.Range("H7").Value = 0.48671
$ws.Range("I7").Value = 0.08672241273057199
$ws.Range("J7").Value = 0.05953597632407945
$ws.Range("M7").Value = 0.4260053333333333
$ws.Range("N7").Value = 1.278016
$ws.Range("O7").Value = 0.02405532912416773
$ws.Range("P7").Value = 0.02531756756689831
$ws.Range("Q7").Value = 0.1036705278933333
$ws.Range("R7").Value = 0.62202316736
$ws.Range("S7").Value = 0.002086136180675822
$ws.Range("T7").Value = 0.001507306103246139

$ws.Range("G8").Value = 0.243355
$ws.Range("H8").Value = 0.48671
$ws.Range("I8").Value = 0.08672241273057199
$ws.Range("J8").Value = 0.05953597632407945
$ws.Range("O8").Value = 0.05879323641880037
$ws.Range("P8").Value = 0.06187825274916518
$ws.Range("Q8").Value = 0.2533794413966666
$ws.Range("R8").Value = 1.52027664838
$ws.Range("S8").Value = 0.005098691314477302
$ws.Range("T8").Value = 0.003683982190649702

$ws.Range("G9").Value = 0.243355
$ws.Range("H9").Value = 0.48671
$ws.Range("I9").Value = 0.08672241273057199
$ws.Range("J9").Value = 0.05953597632407945
$ws.Range("M9").Value = 5.850740666666667
$ws.Range("N9").Value = 17.552222
$ws.Range("O9").Value = 0.3303749538898241
$ws.Range("P9").Value = 0.3477104875323931
$ws.Range("Q9").Value = 1.423806994936667
$ws.Range("R9").Value = 8.54284196962
$ws.Range("S9").Value = 0.02865091310707702
$ws.Range("T9").Value = 0.02070128335336268

$ws.Range("G10").Value = 0.243355
$ws.Range("H10").Value = 0.48671
$ws.Range("I10").Value = 0.08672241273057199
$ws.Range("J10").Value = 0.05953597632407945
$ws.Range("M10").Value = 2.648771
$ws.Range("N10").Value = 5.297542
$ws.Range("O10").Value = 0.1495686865725097
$ws.Range("P10").Value = 0.1049445996947469
$ws.Range("Q10").Value = 0.644591666705
$ws.Range("R10").Value = 2.57836666682
$ws.Range("S10").Value = 0.01297095736851075
$ws.Range("T10").Value = 0.006247979202766444

$ws.Range("G11").Value = 0.243355
$ws.Range("H11").Value = 0.48671
$ws.Range("I11").Value = 0.08672241273057199
$ws.Range("J11").Value = 0.05953597632407945
$ws.Range("M11").Value = 7.742685666666667
$ws.Range("N11").Value = 23.228057
$ws.Range("O11").Value = 0.4372077939946981
$ws.Range("P11").Value = 0.4601490924567965
$ws.Range("Q11").Value = 1.884221270411667
$ws.Range("R11").Value = 11.30532762247
$ws.Range("S11").Value = 0.0379157147598311
$ws.Range("T11").Value = 0.02739542547405448
